$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# En "UsuariosABM" autocompletar campos comunes... -> marcar como resuelto (verde)
$ws.Range("A17").Interior.Color = 65280

# Actualizar el texto de la fila "Averiguar: Como deberia ser el uso del sistema
# para los profesores." y marcarla como IMPORTANTE
$ws.Range("B21").Value = "Averiguar: Cómo deberia ser el uso del sistema para los profesores. A partir de aca seguir con la carga de NOTAS!"
$ws.Range("A21").Value = "IMPORTANTE"

# Actualizar la seleccion/scroll de la hoja
$null = $ws.Range("B23").Select()
